# Daten aktualisiert am 2023-12-24
# Update the cryptocurrency list table (rows 2-51) with refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; Col=4; Value=43605},
    @{Row=2; Col=5; Value=853742003833},
    @{Row=2; Col=6; Value=11057977487},
    @{Row=2; Col=7; Value=-0.12968},
    @{Row=3; Col=4; Value=2289},
    @{Row=3; Col=5; Value=275175557796},
    @{Row=3; Col=6; Value=12143182508},
    @{Row=3; Col=7; Value=-0.23747},
    @{Row=4; Col=4; Value=1},
    @{Row=4; Col=5; Value=91217301562},
    @{Row=4; Col=6; Value=33329855612},
    @{Row=4; Col=7; Value=-0.04523},
    @{Row=5; Col=2; Value="SOL"},
    @{Row=5; Col=3; Value="Solana"},
    @{Row=5; Col=4; Value=109.99},
    @{Row=5; Col=5; Value=46879984480},
    @{Row=5; Col=6; Value=6277323524},
    @{Row=5; Col=7; Value=14.53062},
    @{Row=6; Col=2; Value="BNB"},
    @{Row=6; Col=3; Value="BNB"},
    @{Row=6; Col=4; Value=268.05},
    @{Row=6; Col=5; Value=41186146760},
    @{Row=6; Col=6; Value=564429863},
    @{Row=6; Col=7; Value=-0.466},
    @{Row=7; Col=2; Value="XRP"},
    @{Row=7; Col=3; Value="XRP"},
    @{Row=7; Col=4; Value=0.617929},
    @{Row=7; Col=5; Value=33411209077},
    @{Row=7; Col=6; Value=855018661},
    @{Row=7; Col=7; Value=0.21732},
    @{Row=8; Col=4; Value=0.999838},
    @{Row=8; Col=5; Value=25066662336},
    @{Row=8; Col=6; Value=5984360710},
    @{Row=8; Col=7; Value=-0.04116},
    @{Row=9; Col=4; Value=0.612484},
    @{Row=9; Col=5; Value=21457221904},
    @{Row=9; Col=6; Value=650933583},
    @{Row=9; Col=7; Value=0.84197},
    @{Row=10; Col=4; Value=2283},
    @{Row=10; Col=5; Value=20891958878},
    @{Row=10; Col=6; Value=36533900},
    @{Row=10; Col=7; Value=-0.38421},
    @{Row=11; Col=4; Value=47.3},
    @{Row=11; Col=5; Value=17258054768},
    @{Row=11; Col=6; Value=1636236523},
    @{Row=11; Col=7; Value=4.33921},
    @{Row=12; Col=4; Value=0.094156},
    @{Row=12; Col=5; Value=13407861660},
    @{Row=12; Col=6; Value=594683216},
    @{Row=12; Col=7; Value=1.12547},
    @{Row=13; Col=4; Value=8.98},
    @{Row=13; Col=5; Value=11818629222},
    @{Row=13; Col=6; Value=930229646},
    @{Row=13; Col=7; Value=13.69242},
    @{Row=14; Col=4; Value=0.105826},
    @{Row=14; Col=5; Value=9351114274},
    @{Row=14; Col=6; Value=358350602},
    @{Row=14; Col=7; Value=-0.08069999999999999},
    @{Row=15; Col=4; Value=15.62},
    @{Row=15; Col=5; Value=8701679339},
    @{Row=15; Col=6; Value=688157223},
    @{Row=15; Col=7; Value=1.64231},
    @{Row=16; Col=4; Value=0.842726},
    @{Row=16; Col=5; Value=7821360687},
    @{Row=16; Col=6; Value=484423732},
    @{Row=16; Col=7; Value=-0.17359},
    @{Row=17; Col=2; Value="TON"},
    @{Row=17; Col=3; Value="Toncoin"},
    @{Row=17; Col=4; Value=2.22},
    @{Row=17; Col=5; Value=7680868642},
    @{Row=17; Col=6; Value=34534980},
    @{Row=17; Col=7; Value=-0.80431},
    @{Row=18; Col=2; Value="WBTC"},
    @{Row=18; Col=3; Value="Wrapped Bitcoin"},
    @{Row=18; Col=4; Value=43522},
    @{Row=18; Col=5; Value=6778191618},
    @{Row=18; Col=6; Value=114214627},
    @{Row=18; Col=7; Value=-0.44355},
    @{Row=19; Col=2; Value="SHIB"},
    @{Row=19; Col=3; Value="Shiba Inu"},
    @{Row=19; Col=4; Value=0.00001087},
    @{Row=19; Col=5; Value=6405276409},
    @{Row=19; Col=6; Value=252379549},
    @{Row=19; Col=7; Value=-1.36134},
    @{Row=20; Col=2; Value="LTC"},
    @{Row=20; Col=3; Value="Litecoin"},
    @{Row=20; Col=4; Value=72.16},
    @{Row=20; Col=5; Value=5340922334},
    @{Row=20; Col=6; Value=421855162},
    @{Row=20; Col=7; Value=0.13757},
    @{Row=21; Col=2; Value="DAI"},
    @{Row=21; Col=3; Value="Dai"},
    @{Row=21; Col=4; Value=0.999008},
    @{Row=21; Col=5; Value=5330469213},
    @{Row=21; Col=6; Value=160905099},
    @{Row=21; Col=7; Value=-0.00906},
    @{Row=22; Col=4; Value=6.73},
    @{Row=22; Col=5; Value=5077628403},
    @{Row=22; Col=6; Value=208362209},
    @{Row=22; Col=7; Value=8.45017},
    @{Row=23; Col=4; Value=230.43},
    @{Row=23; Col=5; Value=4512405536},
    @{Row=23; Col=6; Value=140888564},
    @{Row=23; Col=7; Value=-0.52056},
    @{Row=24; Col=2; Value="ICP"},
    @{Row=24; Col=3; Value="Internet Computer"},
    @{Row=24; Col=4; Value=9.92},
    @{Row=24; Col=5; Value=4488023595},
    @{Row=24; Col=6; Value=189680473},
    @{Row=24; Col=7; Value=8.778370000000001},
    @{Row=25; Col=2; Value="NEAR"},
    @{Row=25; Col=3; Value="NEAR Protocol"},
    @{Row=25; Col=4; Value=3.8},
    @{Row=25; Col=5; Value=3828676488},
    @{Row=25; Col=6; Value=621308525},
    @{Row=25; Col=7; Value=13.16078},
    @{Row=26; Col=2; Value="LEO"},
    @{Row=26; Col=3; Value="LEO Token"},
    @{Row=26; Col=4; Value=3.92},
    @{Row=26; Col=5; Value=3633874473},
    @{Row=26; Col=6; Value=1379772},
    @{Row=26; Col=7; Value=0.30423},
    @{Row=27; Col=2; Value="XLM"},
    @{Row=27; Col=3; Value="Stellar"},
    @{Row=27; Col=4; Value=0.126885},
    @{Row=27; Col=5; Value=3590740842},
    @{Row=27; Col=6; Value=54949241},
    @{Row=27; Col=7; Value=1.49727},
    @{Row=28; Col=2; Value="INJ"},
    @{Row=28; Col=3; Value="Injective"},
    @{Row=28; Col=4; Value=41.73},
    @{Row=28; Col=5; Value=3505217737},
    @{Row=28; Col=6; Value=449240568},
    @{Row=28; Col=7; Value=5.27654},
    @{Row=29; Col=2; Value="ATOM"},
    @{Row=29; Col=3; Value="Cosmos Hub"},
    @{Row=29; Col=4; Value=11.56},
    @{Row=29; Col=5; Value=3384782350},
    @{Row=29; Col=6; Value=234914884},
    @{Row=29; Col=7; Value=2.9608},
    @{Row=30; Col=2; Value="OKB"},
    @{Row=30; Col=3; Value="OKB"},
    @{Row=30; Col=4; Value=55.32},
    @{Row=30; Col=5; Value=3316072675},
    @{Row=30; Col=6; Value=7948240},
    @{Row=30; Col=7; Value=-0.91259},
    @{Row=31; Col=2; Value="OP"},
    @{Row=31; Col=3; Value="Optimism"},
    @{Row=31; Col=4; Value=3.51},
    @{Row=31; Col=5; Value=3190796334},
    @{Row=31; Col=6; Value=759367695},
    @{Row=31; Col=7; Value=8.08742},
    @{Row=32; Col=2; Value="XMR"},
    @{Row=32; Col=3; Value="Monero"},
    @{Row=32; Col=4; Value=175.37},
    @{Row=32; Col=5; Value=3186743681},
    @{Row=32; Col=6; Value=87392499},
    @{Row=32; Col=7; Value=-0.18206},
    @{Row=33; Col=4; Value=2.41},
    @{Row=33; Col=5; Value=3119258840},
    @{Row=33; Col=6; Value=136221435},
    @{Row=33; Col=7; Value=-4.6143},
    @{Row=34; Col=2; Value="HBAR"},
    @{Row=34; Col=3; Value="Hedera"},
    @{Row=34; Col=4; Value=0.091975},
    @{Row=34; Col=5; Value=3088396252},
    @{Row=34; Col=6; Value=88096463},
    @{Row=34; Col=7; Value=2.43528},
    @{Row=35; Col=2; Value="ETC"},
    @{Row=35; Col=3; Value="Ethereum Classic"},
    @{Row=35; Col=4; Value=21.06},
    @{Row=35; Col=5; Value=3016100128},
    @{Row=35; Col=6; Value=128156042},
    @{Row=35; Col=7; Value=-1.2539},
    @{Row=36; Col=2; Value="APT"},
    @{Row=36; Col=3; Value="Aptos"},
    @{Row=36; Col=4; Value=9.43},
    @{Row=36; Col=5; Value=2886804660},
    @{Row=36; Col=6; Value=270780693},
    @{Row=36; Col=7; Value=0.93492},
    @{Row=37; Col=2; Value="FIL"},
    @{Row=37; Col=3; Value="Filecoin"},
    @{Row=37; Col=4; Value=5.62},
    @{Row=37; Col=5; Value=2730409211},
    @{Row=37; Col=6; Value=245464098},
    @{Row=37; Col=7; Value=4.14724},
    @{Row=38; Col=2; Value="VET"},
    @{Row=38; Col=3; Value="VeChain"},
    @{Row=38; Col=4; Value=0.03595505},
    @{Row=38; Col=5; Value=2615844748},
    @{Row=38; Col=6; Value=69767280},
    @{Row=38; Col=7; Value=2.77417},
    @{Row=39; Col=2; Value="CRO"},
    @{Row=39; Col=3; Value="Cronos"},
    @{Row=39; Col=4; Value=0.098604},
    @{Row=39; Col=5; Value=2599844959},
    @{Row=39; Col=6; Value=12881767},
    @{Row=39; Col=7; Value=-2.21099},
    @{Row=40; Col=2; Value="TUSD"},
    @{Row=40; Col=3; Value="TrueUSD"},
    @{Row=40; Col=4; Value=0.998974},
    @{Row=40; Col=5; Value=2428705500},
    @{Row=40; Col=6; Value=164430757},
    @{Row=40; Col=7; Value=-0.07778},
    @{Row=41; Col=2; Value="KAS"},
    @{Row=41; Col=3; Value="Kaspa"},
    @{Row=41; Col=4; Value=0.107191},
    @{Row=41; Col=5; Value=2361004532},
    @{Row=41; Col=6; Value=34917686},
    @{Row=41; Col=7; Value=-0.24342},
    @{Row=42; Col=4; Value=2.39},
    @{Row=42; Col=5; Value=2128356383},
    @{Row=42; Col=6; Value=89349215},
    @{Row=42; Col=7; Value=1.91646},
    @{Row=43; Col=2; Value="STX"},
    @{Row=43; Col=3; Value="Stacks"},
    @{Row=43; Col=4; Value=1.41},
    @{Row=43; Col=5; Value=2026192637},
    @{Row=43; Col=6; Value=108071205},
    @{Row=43; Col=7; Value=-3.38912},
    @{Row=44; Col=2; Value="TIA"},
    @{Row=44; Col=3; Value="Celestia"},
    @{Row=44; Col=4; Value=13.71},
    @{Row=44; Col=5; Value=1999770118},
    @{Row=44; Col=6; Value=258674427},
    @{Row=44; Col=7; Value=10.92248},
    @{Row=45; Col=2; Value="MNT"},
    @{Row=45; Col=3; Value="Mantle"},
    @{Row=45; Col=4; Value=0.621167},
    @{Row=45; Col=5; Value=1944487057},
    @{Row=45; Col=6; Value=9482064},
    @{Row=45; Col=7; Value=4.27293},
    @{Row=46; Col=2; Value="ALGO"},
    @{Row=46; Col=3; Value="Algorand"},
    @{Row=46; Col=4; Value=0.238034},
    @{Row=46; Col=5; Value=1906307005},
    @{Row=46; Col=6; Value=108214904},
    @{Row=46; Col=7; Value=0.66772},
    @{Row=47; Col=2; Value="RUNE"},
    @{Row=47; Col=3; Value="THORChain"},
    @{Row=47; Col=4; Value=6.31},
    @{Row=47; Col=5; Value=1894663804},
    @{Row=47; Col=6; Value=814828624},
    @{Row=47; Col=7; Value=22.04921},
    @{Row=48; Col=2; Value="EGLD"},
    @{Row=48; Col=3; Value="MultiversX"},
    @{Row=48; Col=4; Value=71.18000000000001},
    @{Row=48; Col=5; Value=1875088617},
    @{Row=48; Col=6; Value=117866510},
    @{Row=48; Col=7; Value=11.82114},
    @{Row=49; Col=2; Value="FDUSD"},
    @{Row=49; Col=3; Value="First Digital USD"},
    @{Row=49; Col=4; Value=0.998849},
    @{Row=49; Col=5; Value=1793428025},
    @{Row=49; Col=6; Value=1274063858},
    @{Row=49; Col=7; Value=-0.09893},
    @{Row=50; Col=2; Value="RNDR"},
    @{Row=50; Col=3; Value="Render"},
    @{Row=50; Col=4; Value=4.67},
    @{Row=50; Col=5; Value=1749619269},
    @{Row=50; Col=6; Value=203591907},
    @{Row=50; Col=7; Value=6.00432},
    @{Row=51; Col=2; Value="ARB"},
    @{Row=51; Col=3; Value="Arbitrum"},
    @{Row=51; Col=4; Value=1.37},
    @{Row=51; Col=5; Value=1740769212},
    @{Row=51; Col=6; Value=806571831},
    @{Row=51; Col=7; Value=1.82588}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
